$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Donor cells used as format sources for cells whose underlying style
#     (text-placeholder vs integer-count vs percent-change) must switch. ---
$DonorText    = "C15"   # style 13 : General / shared-string placeholder ("0", "***.*")
$DonorInteger = "I16"   # style 14 : "#,##0" integer count
$DonorPercent = "H16"   # style 15 : "#,##0.0;""-""#,##0.0" percent change

function Set-CellWithStyle($cell, $kind, $value, $styleKind) {
    if ($styleKind -eq 13) {
        $ws.Range($DonorText).Copy() | Out-Null
        $ws.Range($cell).PasteSpecial(-4122) | Out-Null
    } elseif ($styleKind -eq 14) {
        $ws.Range($DonorInteger).Copy() | Out-Null
        $ws.Range($cell).PasteSpecial(-4122) | Out-Null
    } elseif ($styleKind -eq 15) {
        $ws.Range($DonorPercent).Copy() | Out-Null
        $ws.Range($cell).PasteSpecial(-4122) | Out-Null
    }

    if ($kind -eq "text") {
        # Force text storage (e.g. the literal string "0") instead of letting
        # Excel auto-coerce a numeric-looking string into a number, then
        # restore the General number format so the stored style matches the
        # other placeholder cells exactly.
        $ws.Range($cell).NumberFormat = "@"
        $ws.Range($cell).Value = [string]$value
        $ws.Range($DonorText).Copy() | Out-Null
        $ws.Range($cell).PasteSpecial(-4122) | Out-Null
    } else {
        $ws.Range($cell).Value = $value
    }
}

$changes = @(
    @{ Cell = "G15"; Kind = "text"; Value = "0"; Style = 13 },
    @{ Cell = "H15"; Kind = "text"; Value = "***.*"; Style = 13 },
    @{ Cell = "L15"; Kind = "number"; Value = 60; Style = 0 },
    @{ Cell = "M15"; Kind = "number"; Value = 0; Style = 0 },
    @{ Cell = "C16"; Kind = "number"; Value = 3; Style = 0 },
    @{ Cell = "D16"; Kind = "number"; Value = 3; Style = 14 },
    @{ Cell = "E16"; Kind = "number"; Value = 0; Style = 15 },
    @{ Cell = "F16"; Kind = "number"; Value = 9; Style = 0 },
    @{ Cell = "G16"; Kind = "number"; Value = 6; Style = 0 },
    @{ Cell = "H16"; Kind = "number"; Value = 50; Style = 0 },
    @{ Cell = "I16"; Kind = "number"; Value = 50; Style = 0 },
    @{ Cell = "J16"; Kind = "number"; Value = 67; Style = 0 },
    @{ Cell = "K16"; Kind = "number"; Value = -25.373134328358; Style = 0 },
    @{ Cell = "L16"; Kind = "number"; Value = -5.660377358490; Style = 0 },
    @{ Cell = "M16"; Kind = "number"; Value = -54.128440366972; Style = 0 },
    @{ Cell = "N16"; Kind = "number"; Value = -85.163204747774; Style = 0 },
    @{ Cell = "C17"; Kind = "number"; Value = 3; Style = 0 },
    @{ Cell = "D17"; Kind = "number"; Value = 3; Style = 0 },
    @{ Cell = "E17"; Kind = "number"; Value = 0; Style = 0 },
    @{ Cell = "F17"; Kind = "number"; Value = 14; Style = 0 },
    @{ Cell = "G17"; Kind = "number"; Value = 12; Style = 0 },
    @{ Cell = "H17"; Kind = "number"; Value = 16.666666666666; Style = 0 },
    @{ Cell = "I17"; Kind = "number"; Value = 71; Style = 0 },
    @{ Cell = "J17"; Kind = "number"; Value = 76; Style = 0 },
    @{ Cell = "K17"; Kind = "number"; Value = -6.578947368421; Style = 0 },
    @{ Cell = "L17"; Kind = "number"; Value = -18.390804597701; Style = 0 },
    @{ Cell = "M17"; Kind = "number"; Value = 14.516129032258; Style = 0 },
    @{ Cell = "N17"; Kind = "number"; Value = -55.063291139240; Style = 0 },
    @{ Cell = "C18"; Kind = "number"; Value = 2; Style = 14 },
    @{ Cell = "D18"; Kind = "text"; Value = "0"; Style = 13 },
    @{ Cell = "E18"; Kind = "text"; Value = "***.*"; Style = 13 },
    @{ Cell = "F18"; Kind = "number"; Value = 3; Style = 0 },
    @{ Cell = "G18"; Kind = "number"; Value = 4; Style = 0 },
    @{ Cell = "H18"; Kind = "number"; Value = -25; Style = 0 },
    @{ Cell = "I18"; Kind = "number"; Value = 34; Style = 0 },
    @{ Cell = "K18"; Kind = "number"; Value = 6.25; Style = 0 },
    @{ Cell = "L18"; Kind = "number"; Value = -32; Style = 0 },
    @{ Cell = "M18"; Kind = "number"; Value = -72.131147540983; Style = 0 },
    @{ Cell = "N18"; Kind = "number"; Value = -94.407894736842; Style = 0 },
    @{ Cell = "C19"; Kind = "number"; Value = 11; Style = 0 },
    @{ Cell = "D19"; Kind = "number"; Value = 11; Style = 0 },
    @{ Cell = "F19"; Kind = "number"; Value = 34; Style = 0 },
    @{ Cell = "H19"; Kind = "number"; Value = -19.047619047619; Style = 0 },
    @{ Cell = "I19"; Kind = "number"; Value = 248; Style = 0 },
    @{ Cell = "J19"; Kind = "number"; Value = 294; Style = 0 },
    @{ Cell = "K19"; Kind = "number"; Value = -15.646258503401; Style = 0 },
    @{ Cell = "L19"; Kind = "number"; Value = -16.778523489932; Style = 0 },
    @{ Cell = "M19"; Kind = "number"; Value = 5.084745762711; Style = 0 },
    @{ Cell = "N19"; Kind = "number"; Value = -25.748502994012; Style = 0 },
    @{ Cell = "C20"; Kind = "number"; Value = 5; Style = 0 },
    @{ Cell = "D20"; Kind = "number"; Value = 2; Style = 0 },
    @{ Cell = "E20"; Kind = "number"; Value = 150; Style = 0 },
    @{ Cell = "F20"; Kind = "number"; Value = 25; Style = 0 },
    @{ Cell = "H20"; Kind = "number"; Value = 78.571428571428; Style = 0 },
    @{ Cell = "I20"; Kind = "number"; Value = 71; Style = 0 },
    @{ Cell = "J20"; Kind = "number"; Value = 76; Style = 0 },
    @{ Cell = "K20"; Kind = "number"; Value = -6.578947368421; Style = 0 },
    @{ Cell = "L20"; Kind = "number"; Value = 36.538461538461; Style = 0 },
    @{ Cell = "M20"; Kind = "number"; Value = -14.457831325301; Style = 0 },
    @{ Cell = "N20"; Kind = "number"; Value = -94.914040114613; Style = 0 },
    @{ Cell = "C21"; Kind = "number"; Value = 24; Style = 0 },
    @{ Cell = "D21"; Kind = "number"; Value = 19; Style = 0 },
    @{ Cell = "E21"; Kind = "number"; Value = 26.315789473684; Style = 0 },
    @{ Cell = "F21"; Kind = "number"; Value = 86; Style = 0 },
    @{ Cell = "G21"; Kind = "number"; Value = 78; Style = 0 },
    @{ Cell = "H21"; Kind = "number"; Value = 10.256410256410; Style = 0 },
    @{ Cell = "I21"; Kind = "number"; Value = 483; Style = 0 },
    @{ Cell = "J21"; Kind = "number"; Value = 550; Style = 0 },
    @{ Cell = "K21"; Kind = "number"; Value = -12.181818181818; Style = 0 },
    @{ Cell = "L21"; Kind = "number"; Value = -11.700182815356; Style = 0 },
    @{ Cell = "M21"; Kind = "number"; Value = -22.096774193548; Style = 0 },
    @{ Cell = "N21"; Kind = "number"; Value = -83.094154707735; Style = 0 },
    @{ Cell = "C23"; Kind = "number"; Value = 2; Style = 14 },
    @{ Cell = "D23"; Kind = "text"; Value = "0"; Style = 13 },
    @{ Cell = "E23"; Kind = "text"; Value = "***.*"; Style = 13 },
    @{ Cell = "F23"; Kind = "number"; Value = 4; Style = 0 },
    @{ Cell = "H23"; Kind = "number"; Value = 33.333333333333; Style = 0 },
    @{ Cell = "I23"; Kind = "number"; Value = 10; Style = 0 },
    @{ Cell = "K23"; Kind = "number"; Value = -50; Style = 0 },
    @{ Cell = "L23"; Kind = "number"; Value = 11.111111111111; Style = 0 },
    @{ Cell = "M23"; Kind = "number"; Value = -41.176470588235; Style = 0 },
    @{ Cell = "C24"; Kind = "number"; Value = 18; Style = 0 },
    @{ Cell = "D24"; Kind = "number"; Value = 55; Style = 0 },
    @{ Cell = "E24"; Kind = "number"; Value = -67.272727272727; Style = 0 },
    @{ Cell = "F24"; Kind = "number"; Value = 81; Style = 0 },
    @{ Cell = "G24"; Kind = "number"; Value = 168; Style = 0 },
    @{ Cell = "H24"; Kind = "number"; Value = -51.785714285714; Style = 0 },
    @{ Cell = "I24"; Kind = "number"; Value = 623; Style = 0 },
    @{ Cell = "J24"; Kind = "number"; Value = 813; Style = 0 },
    @{ Cell = "K24"; Kind = "number"; Value = -23.370233702337; Style = 0 },
    @{ Cell = "L24"; Kind = "number"; Value = 5.414551607445; Style = 0 },
    @{ Cell = "M24"; Kind = "number"; Value = 40.632054176072; Style = 0 },
    @{ Cell = "C25"; Kind = "number"; Value = 13; Style = 0 },
    @{ Cell = "D25"; Kind = "number"; Value = 42; Style = 0 },
    @{ Cell = "E25"; Kind = "number"; Value = -69.047619047619; Style = 0 },
    @{ Cell = "F25"; Kind = "number"; Value = 61; Style = 0 },
    @{ Cell = "G25"; Kind = "number"; Value = 128; Style = 0 },
    @{ Cell = "H25"; Kind = "number"; Value = -52.34375; Style = 0 },
    @{ Cell = "I25"; Kind = "number"; Value = 442; Style = 0 },
    @{ Cell = "J25"; Kind = "number"; Value = 667; Style = 0 },
    @{ Cell = "K25"; Kind = "number"; Value = -33.733133433283; Style = 0 },
    @{ Cell = "L25"; Kind = "number"; Value = 9.950248756218; Style = 0 },
    @{ Cell = "C26"; Kind = "number"; Value = 9; Style = 0 },
    @{ Cell = "D26"; Kind = "number"; Value = 5; Style = 0 },
    @{ Cell = "E26"; Kind = "number"; Value = 80; Style = 0 },
    @{ Cell = "F26"; Kind = "number"; Value = 34; Style = 0 },
    @{ Cell = "G26"; Kind = "number"; Value = 30; Style = 0 },
    @{ Cell = "H26"; Kind = "number"; Value = 13.333333333333; Style = 0 },
    @{ Cell = "I26"; Kind = "number"; Value = 148; Style = 0 },
    @{ Cell = "J26"; Kind = "number"; Value = 160; Style = 0 },
    @{ Cell = "K26"; Kind = "number"; Value = -7.5; Style = 0 },
    @{ Cell = "L26"; Kind = "number"; Value = 13.846153846153; Style = 0 },
    @{ Cell = "M26"; Kind = "number"; Value = -12.941176470588; Style = 0 },
    @{ Cell = "G27"; Kind = "number"; Value = 2; Style = 0 },
    @{ Cell = "L27"; Kind = "number"; Value = -38.461538461538; Style = 0 },
    @{ Cell = "C28"; Kind = "number"; Value = 1; Style = 0 },
    @{ Cell = "I28"; Kind = "number"; Value = 14; Style = 0 },
    @{ Cell = "K28"; Kind = "number"; Value = 7.692307692307; Style = 0 },
    @{ Cell = "L28"; Kind = "number"; Value = -12.5; Style = 0 },
    @{ Cell = "N29"; Kind = "number"; Value = -81.481481481481; Style = 0 },
    @{ Cell = "N30"; Kind = "number"; Value = -79.166666666666; Style = 0 },
    @{ Cell = "G31"; Kind = "number"; Value = 1; Style = 0 },
    @{ Cell = "H31"; Kind = "number"; Value = 0; Style = 0 },
    @{ Cell = "F33"; Kind = "number"; Value = 1; Style = 0 }
)

foreach ($chg in $changes) {
    Set-CellWithStyle $chg.Cell $chg.Kind $chg.Value $chg.Style
}

# --- Header text: volume/issue number and the week-covering date range ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Column E width collapses back to the standard 6.168446 ("best fit")
#     width now that its contents are no longer the widest in the row. ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth()

Write-Output "done"
